# Apply updated cryptocurrency data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.227.13"
$ws.Range("E2").Value = '  -0.08%  '

# Row 3
$ws.Range("D3").Value = "'1.858.59"
$ws.Range("E3").Value = '  -0.28%  '

# Row 4
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
$ws.Range("D5").Value = "'0.7116"
$ws.Range("E5").Value = '  +1.03%  '

# Row 6
$ws.Range("D6").Value = "'241.36"
$ws.Range("E6").Value = '  -0.58%  '

# Row 7
$ws.Range("D7").Value = "'0.9998"
$ws.Range("E7").Value = '  -0.28%  '

# Row 8
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = "'0.3098"
$ws.Range("E8").Value = '  -0.42%  '

# Row 9
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = "'0.07754"
$ws.Range("E9").Value = '  -0.90%  '

# Row 10
$ws.Range("D10").Value = "'23.87"
$ws.Range("E10").Value = '  -1.64%  '

# Row 11
$ws.Range("D11").Value = "'0.07817"
$ws.Range("E11").Value = '  -2.60%  '

# Row 12
$ws.Range("D12").Value = "'1.866.05"
$ws.Range("E12").Value = '  -8.56%  '

# Row 13
$ws.Range("D13").Value = "'5.103"
$ws.Range("E13").Value = '  -1.37%  '

# Row 14
$ws.Range("D14").Value = "'92.05"
$ws.Range("E14").Value = '  -1.05%  '

# Row 15
$ws.Range("D15").Value = "'0.6875"
$ws.Range("E15").Value = '  -1.33%  '

# Row 16
$ws.Range("D16").Value = "'6.506"
$ws.Range("E16").Value = '  +2.56%  '

# Row 17
$ws.Range("D17").Value = "'0.000008417"
$ws.Range("E17").Value = '  +1.62%  '

# Row 18
$ws.Range("D18").Value = "'29.234.25"
$ws.Range("E18").Value = '  -0.90%  '

# Row 19
$ws.Range("D19").Value = "'249.43"
$ws.Range("E19").Value = '  -0.11%  '

# Row 20
$ws.Range("D20").Value = "'2.110.08"
$ws.Range("E20").Value = '  -7.47%  '

# Row 21
$ws.Range("D21").Value = "'12.84"
$ws.Range("E21").Value = '  -2.44%  '

# Row 22
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = '  -0.40%  '

# Row 23
$ws.Range("D23").Value = "'7.515"
$ws.Range("E23").Value = '  -0.91%  '

# Row 24
$ws.Range("D24").Value = "'0.9999"
$ws.Range("E24").Value = '  -0.20%  '

# Row 25
$ws.Range("D25").Value = "'0.1544"
$ws.Range("E25").Value = '  -0.36%  '

# Row 26
$ws.Range("D26").Value = "'160.39"
$ws.Range("E26").Value = '  +0.29%  '

# Row 27
$ws.Range("D27").Value = "'8.859"
$ws.Range("E27").Value = '  -1.32%  '

# Row 28
$ws.Range("D28").Value = "'18.52"
$ws.Range("E28").Value = '  -0.78%  '

# Row 29
$ws.Range("D29").Value = "'1.564"
$ws.Range("E29").Value = '  +4.18%  '

# Row 30
$ws.Range("D30").Value = "'4.249"
$ws.Range("E30").Value = '  -0.65%  '

# Row 31
$ws.Range("D31").Value = "'4.232"
$ws.Range("E31").Value = '  -0.97%  '

# Row 32
$ws.Range("E32").Value = '  -2.14%  '

# Row 33
$ws.Range("D33").Value = "'0.05202"
$ws.Range("E33").Value = '  -0.87%  '

# Row 34
$ws.Range("D34").Value = "'0.7569"
$ws.Range("E34").Value = '  +1.84%  '

# Row 35
$ws.Range("D35").Value = "'1.843"
$ws.Range("E35").Value = '  -2.31%  '

# Row 36
$ws.Range("D36").Value = "'1.164"
$ws.Range("E36").Value = '  +0.49%  '

# Row 37
$ws.Range("D37").Value = "'2.709"
$ws.Range("E37").Value = '  +0.06%  '

# Row 38
$ws.Range("D38").Value = "'0.01858"
$ws.Range("E38").Value = '  +0.09%  '

# Row 39
$ws.Range("D39").Value = "'1.227.02"
$ws.Range("E39").Value = '  -1.96%  '

# Row 40
$ws.Range("E40").Value = '  -0.54%  '

# Row 41
$ws.Range("D41").Value = "'0.8980"
$ws.Range("E41").Value = '  +0.24%  '

# Row 42
$ws.Range("D42").Value = "'109.73"
$ws.Range("E42").Value = '  -1.16%  '

# Row 43
$ws.Range("D43").Value = "'0.9990"
$ws.Range("E43").Value = '  -0.33%  '

# Row 44
$ws.Range("D44").Value = "'5.593"
$ws.Range("E44").Value = '  -10.23%  '

# Row 45
$ws.Range("D45").Value = "'2.006.06"
$ws.Range("E45").Value = '  -4.67%  '

# Row 46
$ws.Range("E46").Value = '  -3.59%  '

# Row 47
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = "'65.10"
$ws.Range("E47").Value = '  -9.87%  '

# Row 48
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = "'0.5179"
$ws.Range("E48").Value = '  -0.39%  '

# Row 49
$ws.Range("D49").Value = "'9.514"
$ws.Range("E49").Value = '  +1.66%  '

# Row 50
$ws.Range("D50").Value = "'1.748"
$ws.Range("E50").Value = '  -2.89%  '

# Row 51
$ws.Range("D51").Value = "'7.003"
$ws.Range("E51").Value = '  +0.39%  '
